$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

$urls = @(
    "https://imgbb.host/images/MS2yG.png",
    "https://imgbb.host/images/M5uVZ.png",
    "https://imgbb.host/images/MSgVS.png",
    "https://imgbb.host/images/MSQ4H.png",
    "https://imgbb.host/images/MSvpp.png",
    "https://imgbb.host/images/MSNe1.png",
    "https://imgbb.host/images/MSH3E.png",
    "https://imgbb.host/images/MSqLc.png",
    "https://imgbb.host/images/MSMgI.png",
    "https://imgbb.host/images/M5wPk.png",
    "https://imgbb.host/images/MSZw6.png"
)

for ($i = 0; $i -lt $urls.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $urls[$i]
}

$ws.Range("B2").Select()
